$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
# Column C header text changes from "Target Max Electricity kWh per anum"
# to "emissions"; column D keeps showing "Star Rating".
$ws.Range("C1").Value = "emissions"
$ws.Range("D1").Value = "Star Rating"

# --- New emissions figures (4 years scope) --------------------------
$ws.Range("C2").Value  = 115563.1
$ws.Range("C3").Value  = 133826.70000000001
$ws.Range("C4").Value  = 152090.4
$ws.Range("C5").Value  = 231126.2
$ws.Range("C6").Value  = 267653.5
$ws.Range("C7").Value  = 304180.8
$ws.Range("C8").Value  = 346689.2
$ws.Range("C9").Value  = 401480.2
$ws.Range("C10").Value = 456271.1
$ws.Range("C11").Value = 462252.3
$ws.Range("C12").Value = 535306.9
$ws.Range("C13").Value = 608361.5
$ws.Range("C14").Value = 577815.4
$ws.Range("C15").Value = 669133.6
$ws.Range("C16").Value = 760451.9

# --- Formatting: drop the old Verdana/#,##0 style on column C and ---
# --- replace it with the default font plus a 2-decimal number format.
$ws.Range("C2:C16").ClearFormats()
$ws.Range("C2:C16").NumberFormat = "#,##0.00"

# --- Match the recorded active selection -----------------------------
$ws.Range("I8").Select()
